# chore: update Sheets via scheduled runner
# Updates leve-profit price/profit figures across the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR sheets (columns H-N: currentAveragePrice*,
# LevePrice*, LeveProfit*) to reflect refreshed market board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 19250
$ws.Range("I21").Value = 19250
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 19250
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -18782

$ws.Range("H23").Value = 19250
$ws.Range("I23").Value = 19250
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 19250
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -19016

$ws.Range("H135").Value = 210
$ws.Range("I135").Value = 210
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 1890
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = 645

$ws.Range("H138").Value = 4114.3335
$ws.Range("I138").Value = 3761.182
$ws.Range("J138").Value = 7999
$ws.Range("K138").Value = 11283.546
$ws.Range("L138").Value = 23997
$ws.Range("M138").Value = -6143.545999999998
$ws.Range("N138").Value = -34277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 7240.7144
$ws.Range("I41").Value = 651.6667
$ws.Range("J41").Value = 12182.5
$ws.Range("K41").Value = 651.6667
$ws.Range("L41").Value = 12182.5
$ws.Range("M41").Value = -237.6667
$ws.Range("N41").Value = -13010.5

$ws.Range("H46").Value = 19500
$ws.Range("I46").Value = 19500
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 19500
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -19181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 1012
$ws.Range("I10").Value = 416
$ws.Range("J10").Value = 2005.3334
$ws.Range("K10").Value = 416
$ws.Range("L10").Value = 2005.3334
$ws.Range("M10").Value = -276
$ws.Range("N10").Value = -2285.3334

$ws.Range("H26").Value = 33593.75
$ws.Range("I26").Value = 33593.75
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 33593.75
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -33301.75
$ws.Range("N26").Value = $null

$ws.Range("H34").Value = 1369
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 1369
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 1369
$ws.Range("N34").Value = -1597

$ws.Range("H36").Value = 5190.1113
$ws.Range("I36").Value = 5470.5
$ws.Range("J36").Value = 2947
$ws.Range("K36").Value = 5470.5
$ws.Range("L36").Value = 2947
$ws.Range("M36").Value = -4936.5
$ws.Range("N36").Value = -4015

$ws.Range("H46").Value = 21999
$ws.Range("I46").Value = 13999
$ws.Range("J46").Value = 29999
$ws.Range("K46").Value = 13999
$ws.Range("L46").Value = 29999
$ws.Range("M46").Value = -13701
$ws.Range("N46").Value = -30595

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 229997.5
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 229997.5
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 229997.5
$ws.Range("N9").Value = -230333.5

$ws.Range("H29").Value = 13991.667
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 13991.667
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 13991.667
$ws.Range("N29").Value = -14577.667

$ws.Range("H31").Value = 3749.75
$ws.Range("I31").Value = 2000
$ws.Range("J31").Value = 8999
$ws.Range("K31").Value = 2000
$ws.Range("L31").Value = 8999
$ws.Range("M31").Value = -1705
$ws.Range("N31").Value = -9589

$ws.Range("H34").Value = 3749.75
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 8999
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 8999
$ws.Range("M34").Value = -1798
$ws.Range("N34").Value = -9403

$ws.Range("H35").Value = 5421.1665
$ws.Range("I35").Value = 5499.8
$ws.Range("J35").Value = 5028
$ws.Range("K35").Value = 5499.8
$ws.Range("L35").Value = 5028
$ws.Range("M35").Value = -5205.8
$ws.Range("N35").Value = -5616

$ws.Range("H38").Value = 20019
$ws.Range("I38").Value = 3358.6667
$ws.Range("J38").Value = 70000
$ws.Range("K38").Value = 3358.6667
$ws.Range("L38").Value = 70000
$ws.Range("M38").Value = -2981.6667
$ws.Range("N38").Value = -70754

$ws.Range("H46").Value = 20019
$ws.Range("I46").Value = 3358.6667
$ws.Range("J46").Value = 70000
$ws.Range("K46").Value = 3358.6667
$ws.Range("L46").Value = 70000
$ws.Range("M46").Value = -3147.6667
$ws.Range("N46").Value = -70422

$ws.Range("H62").Value = 4150
$ws.Range("I62").Value = 4150
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4150
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3526

$ws.Range("H65").Value = 4150
$ws.Range("I65").Value = 4150
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20750
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -17630

$ws.Range("H68").Value = 96166.336
$ws.Range("I68").Value = 90000
$ws.Range("J68").Value = 99249.5
$ws.Range("K68").Value = 90000
$ws.Range("L68").Value = 99249.5
$ws.Range("M68").Value = -89251
$ws.Range("N68").Value = -100747.5

$ws.Range("H71").Value = 96166.336
$ws.Range("I71").Value = 90000
$ws.Range("J71").Value = 99249.5
$ws.Range("K71").Value = 270000
$ws.Range("L71").Value = 297748.5
$ws.Range("M71").Value = -266256
$ws.Range("N71").Value = -305236.5

$ws.Range("H99").Value = 1113533.4
$ws.Range("I99").Value = 1251250
$ws.Range("J99").Value = 1003360
$ws.Range("K99").Value = 1251250
$ws.Range("L99").Value = 1003360
$ws.Range("M99").Value = -1249752
$ws.Range("N99").Value = -1006356

$ws.Range("H126").Value = 1113533.4
$ws.Range("I126").Value = 1251250
$ws.Range("J126").Value = 1003360
$ws.Range("K126").Value = 3753750
$ws.Range("L126").Value = 3010080
$ws.Range("M126").Value = -3751280
$ws.Range("N126").Value = -3015020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 76.5
$ws.Range("I7").Value = 100.5
$ws.Range("J7").Value = 52.5
$ws.Range("K7").Value = 301.5
$ws.Range("L7").Value = 157.5
$ws.Range("M7").Value = -189.5
$ws.Range("N7").Value = -381.5

$ws.Range("H15").Value = 108.71429
$ws.Range("I15").Value = 67.25
$ws.Range("J15").Value = 164
$ws.Range("K15").Value = 201.75
$ws.Range("L15").Value = 492
$ws.Range("M15").Value = -61.75
$ws.Range("N15").Value = -772

$ws.Range("H17").Value = 254.6
$ws.Range("I17").Value = 46.666668
$ws.Range("J17").Value = 566.5
$ws.Range("K17").Value = 140.000004
$ws.Range("L17").Value = 1699.5
$ws.Range("M17").Value = 28.99999600000001
$ws.Range("N17").Value = -2037.5

$ws.Range("H23").Value = 67.23529000000001
$ws.Range("I23").Value = 32.4
$ws.Range("J23").Value = 117
$ws.Range("K23").Value = 97.19999999999999
$ws.Range("L23").Value = 351
$ws.Range("M23").Value = 137.8
$ws.Range("N23").Value = -821

$ws.Range("H26").Value = 418.83334
$ws.Range("I26").Value = 55
$ws.Range("J26").Value = 2238
$ws.Range("K26").Value = 165
$ws.Range("L26").Value = 6714
$ws.Range("M26").Value = 123
$ws.Range("N26").Value = -7290

$ws.Range("H34").Value = 15450
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 15450
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 46350
$ws.Range("N34").Value = -46518

$ws.Range("H46").Value = 3799.1667
$ws.Range("I46").Value = 401.5
$ws.Range("J46").Value = 5498
$ws.Range("K46").Value = 1204.5
$ws.Range("L46").Value = 16494
$ws.Range("M46").Value = -1113.5
$ws.Range("N46").Value = -16676

$ws.Range("H49").Value = 3999
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 3999
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 11997
$ws.Range("N49").Value = -12309

$ws.Range("H50").Value = 405.25
$ws.Range("I50").Value = 334.7143
$ws.Range("J50").Value = 899
$ws.Range("K50").Value = 1004.1429
$ws.Range("L50").Value = 2697
$ws.Range("M50").Value = -523.1428999999999
$ws.Range("N50").Value = -3659

$ws.Range("H53").Value = 405.25
$ws.Range("I53").Value = 334.7143
$ws.Range("J53").Value = 899
$ws.Range("K53").Value = 1004.1429
$ws.Range("L53").Value = 2697
$ws.Range("M53").Value = -523.1428999999999
$ws.Range("N53").Value = -3659

$ws.Range("H75").Value = 4997.5
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 4997.5
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 14992.5
$ws.Range("N75").Value = -16988.5

$ws.Range("H78").Value = 4997.5
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 4997.5
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 44977.5
$ws.Range("N78").Value = -54961.5

$ws.Range("H99").Value = 2435
$ws.Range("I99").Value = 2435
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 7305
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -5059

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 1350
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1350
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1350
$ws.Range("N25").Value = -2408

$ws.Range("H126").Value = 7123.3335
$ws.Range("I126").Value = 7685
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 23055
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -20585
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 5000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 5000
$ws.Range("N24").Value = -5686

$ws.Range("H32").Value = 5782.5
$ws.Range("I32").Value = 3037.1428
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 3037.1428
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -2720.1428
$ws.Range("N32").Value = -25634

$ws.Range("H41").Value = 8000
$ws.Range("I41").Value = 8000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -7562

$ws.Range("H42").Value = 32499.334
$ws.Range("I42").Value = 32498
$ws.Range("J42").Value = 32500
$ws.Range("K42").Value = 32498
$ws.Range("L42").Value = 32500
$ws.Range("M42").Value = -31935
$ws.Range("N42").Value = -33626

$ws.Range("H49").Value = 32499.334
$ws.Range("I49").Value = 32498
$ws.Range("J49").Value = 32500
$ws.Range("K49").Value = 32498
$ws.Range("L49").Value = 32500
$ws.Range("M49").Value = -32351
$ws.Range("N49").Value = -32794

$ws.Range("H55").Value = 627.9091
$ws.Range("I55").Value = 761.6
$ws.Range("J55").Value = 516.5
$ws.Range("K55").Value = 761.6
$ws.Range("L55").Value = 516.5
$ws.Range("M55").Value = -588.6
$ws.Range("N55").Value = -862.5

$ws.Range("H68").Value = 4179.6
$ws.Range("I68").Value = 3349.5
$ws.Range("J68").Value = 4733
$ws.Range("K68").Value = 3349.5
$ws.Range("L68").Value = 4733
$ws.Range("M68").Value = -2600.5
$ws.Range("N68").Value = -6231

$ws.Range("H71").Value = 4179.6
$ws.Range("I71").Value = 3349.5
$ws.Range("J71").Value = 4733
$ws.Range("K71").Value = 16747.5
$ws.Range("L71").Value = 23665
$ws.Range("M71").Value = -13003.5
$ws.Range("N71").Value = -31153

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 33249.25
$ws.Range("I124").Value = 32999
$ws.Range("J124").Value = 33499.5
$ws.Range("K124").Value = 32999
$ws.Range("L124").Value = 33499.5
$ws.Range("M124").Value = -28089
$ws.Range("N124").Value = -43319.5

$ws.Range("H126").Value = 4640
$ws.Range("I126").Value = 5795
$ws.Range("J126").Value = 3485
$ws.Range("K126").Value = 17385
$ws.Range("L126").Value = 10455
$ws.Range("M126").Value = -14915
$ws.Range("N126").Value = -15395
